# Regenerate the "K" column (strikeouts, column G) values for the
# banda_anthony 2023 save_data sheet.
#
# The save_data regeneration process now sources strikeout counts ("K")
# from the official play-by-play data instead of estimating Strike# from
# pitch-count heuristics, so column G (header "K") is rewritten with the
# recalculated values for every data row (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 3
    16 = 1
    17 = 3
    18 = 0
    19 = 1
    20 = 2
    21 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
